$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 2-19 hold fastqFileName values in column F. All but the one that
# already reads "..._sequence.txt_XXXXXXX.fq" (row 6) need ".txt" inserted
# right before the trailing "_<index>.fq" suffix so they match /lts
# (i.e. "..._sequence_XXXXXXX.fq" -> "..._sequence.txt_XXXXXXX.fq").
for ($row = 2; $row -le 19; $row++) {
    $cell = $ws.Cells.Item($row, 6)
    $current = [string]$cell.Value2
    if ($current -notlike "*.txt_*") {
        $updated = $current -replace "_sequence_", "_sequence.txt_"
        $cell.Value2 = $updated
    }
}

# Row heights for rows 6 and 8 shrink from 15 to 13.8.
$ws.Rows.Item(6).RowHeight = 13.8
$ws.Rows.Item(8).RowHeight = 13.8

# Move the active selection from F19 to H30.
$ws.Range("H30").Select()
